# Site updated: 2025年10月8日 23:49:37
#
# The "1-1疯狂年代" transcript sheet gets two of its Chinese (translated)
# subtitle cells (F30 and F31) expanded: the original machine translation
# line is kept and a revised/alternate translation line is appended below
# it (separated by a line break), and the cells are switched to a
# word-wrapping style so both lines are visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# --- F30: "对于红卫兵来说..." paragraph gets a second, revised translation line ---
$f30 = "对于红卫兵来说，在后两个心理阶段对受害者进行大量虐待是非常无聊的。只有那些仍处于初级阶段的怪物和恶魔才能给他们过度刺激的大脑带来他们渴望的刺激，就像斗牛士的红色斗篷。" + $nl + `
        "对于红卫兵来说，对处于后两个精神阶段的受害者进行虐待是非常无聊的。只有那些仍处于初始阶段的妖魔鬼怪才能给他们过度兴奋的大脑带来他们渴望的刺激，就像斗牛士的红色斗篷一样。" + $nl

$ws.Range("F30").Value = $f30
$ws.Range("F30").WrapText = $true

# --- F31: "但是这种理想的受害者..." paragraph gets a second, revised translation line ---
$f31 = "但是这种理想的受害者越来越少了。清华大概只剩下一个了。因为他是如此罕见，他被保留到最后的斗争会议。" + $nl + `
        "但这种理想的牺牲品已经越来越少了。在清华，可能只剩下一个了。因为他太稀有了，所以被保留到了斗争会的最后。" + $nl

$ws.Range("F31").Value = $f31
$ws.Range("F31").WrapText = $true

# --- View state: the author's selection ended up on D31 after the edit ---
[void]$ws.Range("D31").Select()
